# Update the Metadata worksheet with the latest IG publisher run values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental flag (row 7, next to "Experimental") is now populated with the
# literal text "false". A leading apostrophe forces Excel to store this as
# text instead of auto-converting it to a Boolean, and re-pasting the
# formatting from a neighboring cell keeps the original cell style intact.
$expCell = $ws.Range("B7")
$expCell.Value = "'false"
$ws.Range("B6").Copy() | Out-Null
$expCell.PasteSpecial(-4122) | Out-Null

# Date of this IG publication run (row 8, next to "Date") was refreshed
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# Description (row 12) now carries the IG title text
$ws.Range("B12").Value = "Type of Condition"

$excel.CutCopyMode = 0
